$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.674501419067383
$ws.Range("B1").Value = 2.353157520294189
$ws.Range("C1").Value = 3.119368314743042
$ws.Range("D1").Value = 3.722827672958374
$ws.Range("E1").Value = 0.5787047743797302
